$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at P (between SEX and SHELL) and add the new
# EGG_STATE header, shifting SHELL..COMMENTS one column to the right.
$ws.Columns("P").Insert()
$ws.Cells.Item(1, 16).Value = "EGG_STATE"

# Give the new EGG_STATE column a custom width.
$ws.Columns("P").ColumnWidth = 10.45

# Match the author's final cursor position/selection.
[void]$ws.Range("T8").Select()

# Add the "Unclassified" classification banner to the right header.
$ws.PageSetup.RightHeader = "&""Calibri""&12&K000000 Unclassified - Non-Classifié&1#" + [char]0x0D
